$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista")

$ws.Range("A8").Value = "Empresa 3"
$ws.Range("B8").Value = "Papitas"
$ws.Range("C8").Value = 200

$ws.Range("D8").Select()
